# Adds three new sheets to the workbook:
#   CypherOutput_Message  - identical "run log" layout to the existing Message sheet
#   StatOutput            - a small 4-column stats table (file/sample/case/study counts)
#   StatOutput_Message    - same "run log" layout, but for the stats Cypher query, doubled up
#
# This mirrors the existing "Message" sheet content/structure.

$wb = $excel.ActiveWorkbook

# ---- text shared by the "run log" style sheets (Message / CypherOutput_Message / StatOutput_Message) ----
$neo4jUrlLabel   = "Neo4j_URL:"
$neo4jUrlValue   = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$userNameLabel   = "User_name:"
$userNameValue   = "neo4j"
$pwdLabel        = "PWD:"
$pwdValue        = "icdcDBneo4j0"
$cypherLabel     = "Cypher:"
$cypherQuery     = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN ['Rottweiler'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(s.clinical_study_designation,'') AS ``Study Code`` , coalesce(s.clinical_study_type,'') AS  ``Study Type``, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS ``Stage of Disease`` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  ``Neutered Status``"
$outputLabel     = "Output:"
$outputValue     = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC33_Canine_Filter_Breed-Rottweiler_Neo4jData.xlsx"

$statCypherQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Rottweiler']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

function Write-MessageSheet($ws, [int]$startRow) {
    $ws.Cells.Item($startRow,     1).Value = $neo4jUrlLabel
    $ws.Cells.Item($startRow + 1, 1).Value = $neo4jUrlValue
    $ws.Cells.Item($startRow + 2, 1).Value = $userNameLabel
    $ws.Cells.Item($startRow + 3, 1).Value = $userNameValue
    $ws.Cells.Item($startRow + 4, 1).Value = $pwdLabel
    $ws.Cells.Item($startRow + 5, 1).Value = $pwdValue
    $ws.Cells.Item($startRow + 6, 1).Value = $cypherLabel
    $ws.Cells.Item($startRow + 7, 1).Value = $cypherQuery
    $ws.Cells.Item($startRow + 8, 1).Value = $outputLabel
    $ws.Cells.Item($startRow + 9, 1).Value = $outputValue
}

# ---- 1. CypherOutput_Message : same 10-row "run log" as Message ----
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cypherMsgSheet = $wb.Worksheets.Add($null, $afterSheet)
$cypherMsgSheet.Name = "CypherOutput_Message"
Write-MessageSheet $cypherMsgSheet 1

# ---- 2. StatOutput : 2-row x 4-col stats table (counts stored as text, like the source data) ----
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statSheet = $wb.Worksheets.Add($null, $afterSheet)
$statSheet.Name = "StatOutput"

$statHeaderRange = $statSheet.Range("A1:D2")
$statHeaderRange.NumberFormat = "@"

$statSheet.Cells.Item(1, 1).Value = "number_of_files"
$statSheet.Cells.Item(1, 2).Value = "number_of_sample"
$statSheet.Cells.Item(1, 3).Value = "number_of_cases"
$statSheet.Cells.Item(1, 4).Value = "number_of_study"
$statSheet.Cells.Item(2, 1).Value = "20"
$statSheet.Cells.Item(2, 2).Value = "4"
$statSheet.Cells.Item(2, 3).Value = "6"
$statSheet.Cells.Item(2, 4).Value = "2"

# drop the text-forcing number format now that the values are committed as text
$statHeaderRange.ClearFormats()

# ---- 3. StatOutput_Message : run log (rows 1-10) + run log for the stats query (rows 11-20) ----
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statMsgSheet = $wb.Worksheets.Add($null, $afterSheet)
$statMsgSheet.Name = "StatOutput_Message"

Write-MessageSheet $statMsgSheet 1

$statMsgSheet.Cells.Item(11, 1).Value = $neo4jUrlLabel
$statMsgSheet.Cells.Item(12, 1).Value = $neo4jUrlValue
$statMsgSheet.Cells.Item(13, 1).Value = $userNameLabel
$statMsgSheet.Cells.Item(14, 1).Value = $userNameValue
$statMsgSheet.Cells.Item(15, 1).Value = $pwdLabel
$statMsgSheet.Cells.Item(16, 1).Value = $pwdValue
$statMsgSheet.Cells.Item(17, 1).Value = $cypherLabel
$statMsgSheet.Cells.Item(18, 1).Value = $statCypherQuery
$statMsgSheet.Cells.Item(19, 1).Value = $outputLabel
$statMsgSheet.Cells.Item(20, 1).Value = $outputValue

# restore the original active/selected sheet (CypherOutput) so the new sheets don't steal focus
$wb.Worksheets.Item("CypherOutput").Activate()
